# Update values that were corrected in the source table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: 4 -> 8
$ws.Range("A2").Value = 8

# B3: 7 -> 9
$ws.Range("B3").Value = 9

# Remove the now-empty trailing rows 4 and 5 (shrinks used range / dimension to A1:E3)
$ws.Rows("4:5").Delete() | Out-Null

# Update the active selection shown when the workbook is reopened
$ws.Range("B5").Select() | Out-Null
